# Update the "quiz" marksheet's Correct/Total marking row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Marking" row (row 11): right-answer score weight 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): total score 54 -> 90, and the Corr/Total label 51/84 -> 90/140
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
